$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold + border + centered) from the last existing
# header cell (AC1) onto the three new header cells so they match the rest
# of row 1's formatting.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record columns: every player row gets the team's overall record
$ws.Range("AD2:AD56").Value = 103
$ws.Range("AE2:AE56").Value = 59
$ws.Range("AF2:AF56").Value = 0
